# #5: insurance, claim, debt, investment done
#
# This script finishes the "insurance" (保險) sheet so that it follows the
# same normalized schema as the other property sheets in the workbook
# (company/name/owner/property_category/category/date/legislator_name/
# legislator_id/source_file/index), and fixes the "具有相當價值之財產"
# (other valuables) sheet's category label from the placeholder
# "otherbonds" to "antique".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "保險" (insurance) - 7th sheet
# ---------------------------------------------------------------------
$ins = $wb.Worksheets.Item(7)

# --- Header row -------------------------------------------------------
# Correct the mislabeled header cells (B1:E1) and extend the header with
# the standard metadata columns used throughout the workbook (F1:K1),
# copying the existing header formatting so the new cells share style 1.
$ins.Range("B1:E1").Copy()
$ins.Range("F1:K1").PasteSpecial(-4122)

$ins.Range("B1").Value = "company"
$ins.Range("C1").Value = "name"
$ins.Range("D1").Value = "owner"
$ins.Range("E1").Value = "property_category"
$ins.Range("F1").Value = "category"
$ins.Range("G1").Value = "date"
$ins.Range("H1").Value = "legislator_name"
$ins.Range("I1").Value = "legislator_id"
$ins.Range("J1").Value = "source_file"
$ins.Range("K1").Value = "index"

# --- Data rows (2-9) ---------------------------------------------------
# Column E held ad-hoc free-text insurance policy details; normalize it to
# the shared "insurance" category value (matching the other sheets), and
# append the standard metadata columns F:K.
$insIndex = @{2=113; 3=114; 4=115; 5=116; 6=117; 7=118; 8=120; 9=121}

foreach ($r in 2..9) {
    $ins.Range("B$r`:E$r").Copy()
    $ins.Range("F$r`:K$r").PasteSpecial(-4122)

    $ins.Range("E$r").Value = "insurance"
    $ins.Range("F$r").Value = "normal"

    # "2012-03-03" must stay a plain text value (as it is on every other
    # sheet) rather than being auto-converted into a date serial number:
    # force text format while assigning it, then restore the cell's
    # formatting back to the plain (unformatted) state used by its
    # neighbours.
    $ins.Range("G$r").NumberFormat = "@"
    $ins.Range("G$r").Value = "2012-03-03"
    $ins.Range("F$r").Copy()
    $ins.Range("G$r").PasteSpecial(-4122)

    $ins.Range("H$r").Value = "孫大千"
    $ins.Range("I$r").Value = 919
    $ins.Range("J$r").Value = "tmpc261"
    $ins.Range("K$r").Value = $insIndex[$r]
}

# ---------------------------------------------------------------------
# Sheet "具有相當價值之財產" (other valuables) - 6th sheet
# ---------------------------------------------------------------------
# The "category" column (F) used the placeholder label "otherbonds";
# rename it to "antique" for every data row.
$val = $wb.Worksheets.Item(6)
foreach ($r in 2..6) {
    $val.Range("F$r").Value = "antique"
}
